# Auto-update hourly job matches and history: 2026-02-21 09:39

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-7 (Title, Company, Location, Match Score, Matched Keywords, Posted At, Apply Link)
$rows = @(
    @("Ai/ML Engineer", "Johnson Controls", "Dallas, TX, US USA", 20, "AI Engineer, Data Scientist, Generative AI, LangChain, RAG, LLaMA, FAISS, PyTorch, Azure ML, Docker", "2026-02-21", "https://www.indeed.com/viewjob?jk=1a93129ac0f94b09"),
    @("Backend Software Engineer (hybrid)", "Johnson Controls", "Milwaukee, WI, US USA", 15.6, "RAG, Docker, Kubernetes, CI/CD, Jenkins, GitHub Actions, Git, PostgreSQL, MySQL, MongoDB", "2026-02-21", "https://www.indeed.com/viewjob?jk=f5227f1e48c9e496"),
    @("Data Scientist - Kaggle Grandmaster", "YO IT CONSULTING", "Remote, US USA", 12.2, "Data Scientist, BigQuery, Snowflake, BigQuery, Polars, Python, SQL, R, Scala, Bayesian", "2026-02-21", "https://www.indeed.com/viewjob?jk=a7c00bc5a4960085"),
    @("Application Development Intern - Artificial Intelligence", "C1", "Remote, US USA", 11.1, "Copilot, TensorFlow, Keras, NLTK, Git, Kafka, MongoDB, Python, R, Java", "2026-02-20", "https://www.indeed.com/viewjob?jk=572965e66b126cf4"),
    @("Software Engineer", "Ascension", "Remote, US USA", 11.1, "RAG, CI/CD, Jenkins, Git, MongoDB, NoSQL, SQL, R, Java, Scala", "2026-02-20", "https://www.indeed.com/viewjob?jk=11c257b267f7ca14"),
    @("Application Development Intern - Artificial Intelligence", "C1", "Remote, US USA", 10, "TensorFlow, Keras, NLTK, Git, Kafka, MongoDB, Python, R, Java", "2026-02-20", "https://www.indeed.com/viewjob?jk=b3168e2727d0f09d")
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    # Posted At is stored as plain text (e.g. "2026-02-21"), not a date value.
    # Force text typing via a temporary "@" number format so the date-like
    # string isn't auto-converted into a date serial, then restore the
    # default "Normal" style so no stray style index is left on the cell.
    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value = $rowData[5]
    $ws.Cells.Item($r, 6).Style = "Normal"
    $ws.Cells.Item($r, 7).Value = $rowData[6]
}
